# Included MTM derived from orders & option chain prices.
#
# This script:
#  - updates the 4 existing open-position rows on "portfolio_spy1016dls"
#    (row 2-5) so the timestamp columns (W/X/Y) point at the ORDER
#    execution time (9/19) rather than the later mark time (9/20), and
#    fixes up the stale mark-to-market (J) values,
#  - appends 4 new rows (6-9) that hold the ORIGINAL 9/20 mark-to-market
#    snapshot (what rows 2-5 used to contain) keyed off the option chain,
#  - adds a derived "mtm" column Z = L*C (price * avgCost-ish multiplier)
#    for every data row,
#  - tidies up the various sheet selections / active sheet that Excel
#    persists, and sets a page setup on the orders/portfolio sheets.

$wb = $excel.ActiveWorkbook

$summary   = $wb.Worksheets.Item("summary")
$account   = $wb.Worksheets.Item("account_spy1016dls")
$orders    = $wb.Worksheets.Item("orders_spy1016dls")
$portfolio = $wb.Worksheets.Item("portfolio_spy1016dls")

# ---------------------------------------------------------------------
# portfolio_spy1016dls : fix existing rows 2-5
# ---------------------------------------------------------------------

# New mark-to-market values for the still-open rows.
$portfolio.Range("J2").Value = -550.4
$portfolio.Range("J3").Value = 95.47
$portfolio.Range("J4").Value = 1159.81
$portfolio.Range("J5").Value = -2059.71

# W/X/Y used to carry the 9/20 15:42:02 snapshot timestamp (shared with
# the new rows below); rows 2-5 now carry the 9/19 21:42:02 order time.
for ($r = 2; $r -le 5; $r++) {
    $portfolio.Range("W$r").Value = 20160919
    $portfolio.Range("X$r").Value = 20160919214202
    $portfolio.Range("X$r").NumberFormat = "0"
    $portfolio.Range("Y$r").Value = 42632.903541666703
}

# ---------------------------------------------------------------------
# portfolio_spy1016dls : append rows 6-9 (the old 9/20 snapshot, now
# kept as its own set of rows instead of overwriting rows 2-5)
# ---------------------------------------------------------------------

$newRows = @(
    @{ Row=6; Src=2; A="=A2+100"; B="DU242089"; C="14.21597335"; F="237606530"; G="20161021"; H="SPY   161021C00223000"; I="0.1868009"; J=-560.4;    K="100"; L="-30"; M="AMEX"; N="0.0"; O="C"; R="OPT"; S="223.0"; T="SPY"; V="-133.92" },
    @{ Row=7; Src=3; A="=A3+100"; B="DU242089"; C="4.3317";      F="237606550"; G="20161021"; H="SPY   161021C00228000"; I="0.0301575"; J=90.47;   K="100"; L="30";  M="AMEX"; N="0.0"; O="C"; R="OPT"; S="228.0"; T="SPY"; V="-39.48" },
    @{ Row=8; Src=4; A="=A4+100"; B="DU242089"; C="44.3117";     F="237606986"; G="20161021"; H="SPY   161021P00195000"; I="0.39993785"; J=1199.81; K="100"; L="30";  M="AMEX"; N="0.0"; O="P"; R="OPT"; S="195.0"; T="SPY"; V="-129.54" },
    @{ Row=9; Src=5; A="=A5+100"; B="DU242089"; C="72.48471";    F="237607006"; G="20161021"; H="SPY   161021P00200000"; I="0.66990495"; J=-2009.71; K="100"; L="-30"; M="AMEX"; N="0.0"; O="P"; R="OPT"; S="200.0"; T="SPY"; V="164.83" }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $portfolio.Range("A$r").Formula = $row.A
    $portfolio.Range("A$r").NumberFormat = "General"
    $portfolio.Range("B$r").Value = $row.B
    $portfolio.Range("C$r").Value = $row.C
    $portfolio.Range("E$r").Value = "None"
    $portfolio.Range("F$r").Value = $row.F
    $portfolio.Range("G$r").Value = $row.G
    $portfolio.Range("H$r").Value = $row.H
    $portfolio.Range("I$r").Value = $row.I
    $portfolio.Range("J$r").Value = $row.J
    $portfolio.Range("K$r").Value = $row.K
    $portfolio.Range("L$r").Value = $row.L
    $portfolio.Range("M$r").Value = $row.M
    $portfolio.Range("N$r").Value = $row.N
    $portfolio.Range("O$r").Value = $row.O
    $portfolio.Range("R$r").Value = $row.R
    $portfolio.Range("S$r").Value = $row.S
    $portfolio.Range("T$r").Value = $row.T
    $portfolio.Range("U$r").Value = "None"
    $portfolio.Range("V$r").Value = $row.V
    $portfolio.Range("W$r").Value = 20160920
    $portfolio.Range("X$r").Value = 20160920154202
    $portfolio.Range("X$r").NumberFormat = "0"
    $portfolio.Range("Y$r").Value = 42633.654166666704
}

# ---------------------------------------------------------------------
# portfolio_spy1016dls : derived mtm column Z = L * C for every row
# ---------------------------------------------------------------------
$portfolio.Range("Z2:Z9").Formula = "=L2*C2"

# Column X widened to fit the new numeric timestamps, column layout grew
$portfolio.Columns.Item("X").ColumnWidth = 17.6640625

# ---------------------------------------------------------------------
# Sheet views / selections the commit also recorded
# ---------------------------------------------------------------------
$account.Range("F10").Select()

$orders.PageSetup.PaperSize = 9
$orders.PageSetup.Orientation = 1
$orders.Range("A7").Select()

$portfolio.PageSetup.PaperSize = 9
$portfolio.PageSetup.Orientation = 1
$portfolio.Range("J6").Select()

# Make portfolio_spy1016dls the active sheet/tab (was orders_spy1016dls).
$portfolio.Activate()
